# Updating and adding data translators and more testing
# - bump the "maxiter" row (row 8, columns B:G) from 20 to 2000
# - move the active selection to E17

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B8:G8").Value = 2000

$ws.Range("E17").Select()
